# Add four "PE n" labels to slide 1, matching the lab4_fig diagram update.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$labels = @(
    @{ Text = "PE 0"; Left = 2093874 / 12700; Top = 2719732 / 12700; Name = "TextBox 89" },
    @{ Text = "PE 1"; Left = 3547207 / 12700; Top = 2722001 / 12700; Name = "TextBox 90" },
    @{ Text = "PE 2"; Left = 2084653 / 12700; Top = 1251985 / 12700; Name = "TextBox 91" },
    @{ Text = "PE 3"; Left = 3561294 / 12700; Top = 1260134 / 12700; Name = "TextBox 92" }
)

$width = 531274 / 12700
$height = 307777 / 12700

foreach ($label in $labels) {
    $box = $s.Shapes.AddTextbox(1, $label.Left, $label.Top, $width, $height)
    $box.Name = $label.Name
    $box.Fill.Visible = $false
    $box.TextFrame.WordWrap = $true
    $box.TextFrame.AutoSize = 1
    $tr = $box.TextFrame.TextRange
    $tr.Text = $label.Text
    $tr.ParagraphFormat.Alignment = 2
    $tr.Font.Size = 14
}
